$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update grand final tip dates
$ws.Range("A3").Value = 45920
$ws.Range("A5").Value = 45920
$ws.Range("A9").Value = 45919
$ws.Range("A11").Value = 45919
